$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (rows 1-7, columns A-D)
$data = @(
    @("1.000 GHz", "21.297 dB", "38803.457K", "0.681 dB"),
    @("1.500 GHz", "21.785 dB", "43451.250 K", "0.654 dB"),
    @("2.000 GHz", "22.051 dB", "46219.216 K", "1.201 dB"),
    @("2.500 GHz", "21.564 dB", "41284.808 K", "1.981 dB"),
    @("3.000 GHz", "22.348 dB", "49511.121 K", "2.535 dB"),
    @("3.500 GHz", "22.143 dB", "47208.777 K", "3.382 dB"),
    @("4.000 GHz", "21.559 dB", "41234.481 K", "2.742 dB")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Adjust column C width to match new best-fit width after content change
# (content got shorter, so Excel's automatic best-fit shrinks the column)
$ws.Columns.Item(3).ColumnWidth = 10.33
